$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028588243299612
$ws.Range("D2").Value = 1.03202567865613
$ws.Range("E2").Value = 1.028535837031138
$ws.Range("I2").Value = 1.033665697739855
$ws.Range("J2").Value = 1.033739728823548
$ws.Range("K2").Value = 1.034832099464689
$ws.Range("L2").Value = 1.031352369682536
$ws.Range("N2").Value = 1.015255433175407
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029442182218056
$ws.Range("D3").Value = 1.032640011072422
$ws.Range("E3").Value = 1.029258410206078
$ws.Range("I3").Value = 1.033838246318672
$ws.Range("J3").Value = 1.034234829934326
$ws.Range("K3").Value = 1.035255638743594
$ws.Range("L3").Value = 1.031883129228256
$ws.Range("N3").Value = 1.015420450910096
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029995186491516
$ws.Range("D4").Value = 1.03303786940265
$ws.Range("E4").Value = 1.029726739329396
$ws.Range("I4").Value = 1.033948912898976
$ws.Range("J4").Value = 1.034554998357283
$ws.Range("K4").Value = 1.035529348177799
$ws.Range("L4").Value = 1.032226672669674
$ws.Range("N4").Value = 1.015527128435241
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030227775376865
$ws.Range("D5").Value = 1.033205209743408
$ws.Range("E5").Value = 1.029923808907969
$ws.Range("I5").Value = 1.033995200912261
$ws.Range("J5").Value = 1.034689549196025
$ws.Range("K5").Value = 1.035644330914744
$ws.Range("L5").Value = 1.032371122491462
$ws.Range("N5").Value = 1.015571951271209
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030266834239132
$ws.Range("D6").Value = 1.033233311605742
$ws.Range("E6").Value = 1.029956908497128
$ws.Range("I6").Value = 1.034002959001967
$ws.Range("J6").Value = 1.034712138024549
$ws.Range("K6").Value = 1.035663632019774
$ws.Range("L6").Value = 1.03239537762978
$ws.Range("N6").Value = 1.015579475778665
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029998293941426
$ws.Range("D7").Value = 1.033040105097521
$ws.Range("E7").Value = 1.0297293718619
$ws.Range("I7").Value = 1.033949532330685
$ws.Range("J7").Value = 1.034556796422229
$ws.Range("K7").Value = 1.03553088491712
$ws.Range("L7").Value = 1.032228602722022
$ws.Range("N7").Value = 1.015527727456604
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02887674213054
$ws.Range("D8").Value = 1.032233222888044
$ws.Range("E8").Value = 1.028779872297904
$ws.Range("I8").Value = 1.03372421465593
$ws.Range("J8").Value = 1.033907090202616
$ws.Range("K8").Value = 1.034975308043718
$ws.Range("L8").Value = 1.031531719492103
$ws.Range("N8").Value = 1.015311222117693
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026903928611014
$ws.Range("D9").Value = 1.030814104315396
$ws.Range("E9").Value = 1.027112753080928
$ws.Range("I9").Value = 1.033319670464405
$ws.Range("J9").Value = 1.032760786184878
$ws.Range("K9").Value = 1.033993692476339
$ws.Range("L9").Value = 1.030304597388618
$ws.Range("N9").Value = 1.014928966573706
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025591166220983
$ws.Range("D10").Value = 1.029869949069431
$ws.Range("E10").Value = 1.026005492303557
$ws.Range("I10").Value = 1.033044967998966
$ws.Range("J10").Value = 1.031995691037698
$ws.Range("K10").Value = 1.033337593524155
$ws.Range("L10").Value = 1.02948718430143
$ws.Range("N10").Value = 1.014673655617424
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02502332488399
$ws.Range("D11").Value = 1.029461597055437
$ws.Range("E11").Value = 1.025527042969119
$ws.Range("I11").Value = 1.032924839771834
$ws.Range("J11").Value = 1.031664200427896
$ws.Range("K11").Value = 1.033053110014338
$ws.Range("L11").Value = 1.029133410740969
$ws.Range("N11").Value = 1.014562996482377
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024812493993153
$ws.Range("D12").Value = 1.02930998978737
$ws.Range("E12").Value = 1.025349477885327
$ws.Range("I12").Value = 1.032880042140984
$ws.Range("J12").Value = 1.031541041358925
$ws.Range("K12").Value = 1.03294738303513
$ws.Range("L12").Value = 1.029002030635086
$ws.Range("N12").Value = 1.014521877014438
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024857713807917
$ws.Range("D13").Value = 1.029342506736684
$ws.Range("E13").Value = 1.025387559278027
$ws.Range("I13").Value = 1.032889659368431
$ws.Range("J13").Value = 1.031567460675739
$ws.Range("K13").Value = 1.032970064397289
$ws.Range("L13").Value = 1.029030210863686
$ws.Range("N13").Value = 1.014530697984718
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025005895678522
$ws.Range("D14").Value = 1.029449063655257
$ws.Range("E14").Value = 1.025512362266704
$ws.Range("I14").Value = 1.032921140385652
$ws.Range("J14").Value = 1.031654020638986
$ws.Range("K14").Value = 1.033044371747863
$ws.Range("L14").Value = 1.029122550258853
$ws.Range("N14").Value = 1.014559597852732
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025097207478728
$ws.Range("D15").Value = 1.029514726594668
$ws.Range("E15").Value = 1.025589277741724
$ws.Range("I15").Value = 1.03294051349016
$ws.Range("J15").Value = 1.031707349324106
$ws.Range("K15").Value = 1.033090147423847
$ws.Range("L15").Value = 1.02917944724531
$ws.Range("N15").Value = 1.01457740194349
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025628864612069
$ws.Range("D16").Value = 1.029897060169673
$ws.Range("E16").Value = 1.026037266670745
$ws.Range("I16").Value = 1.033052915703644
$ws.Range("J16").Value = 1.032017686884533
$ws.Range("K16").Value = 1.033356465678965
$ws.Range("L16").Value = 1.029510666820318
$ws.Range("N16").Value = 1.014680997467573
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025962519029555
$ws.Range("D17").Value = 1.030137015890328
$ws.Range("E17").Value = 1.026318547781247
$ws.Range("I17").Value = 1.03312310712255
$ws.Range("J17").Value = 1.032212300948301
$ws.Range("K17").Value = 1.03352341696774
$ws.Range("L17").Value = 1.029718479031562
$ws.Range("N17").Value = 1.014745951663981
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026157191163502
$ws.Range("D18").Value = 1.030277023566978
$ws.Range("E18").Value = 1.026482710842654
$ws.Range("I18").Value = 1.033163934668606
$ws.Range("J18").Value = 1.032325796719373
$ws.Range("K18").Value = 1.033620759280478
$ws.Range("L18").Value = 1.029839708854589
$ws.Range("N18").Value = 1.014783827915352
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026223578979445
$ws.Range("D19").Value = 1.030324770242656
$ws.Range("E19").Value = 1.026538702526377
$ws.Range("I19").Value = 1.033177836463681
$ws.Range("J19").Value = 1.032364492529171
$ws.Range("K19").Value = 1.033653944087845
$ws.Range("L19").Value = 1.029881047844923
$ws.Range("N19").Value = 1.014796740954774
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025926715145308
$ws.Range("D20").Value = 1.03011126618703
$ws.Range("E20").Value = 1.026288358980795
$ws.Range("I20").Value = 1.033115588026213
$ws.Range("J20").Value = 1.032191422690335
$ws.Range("K20").Value = 1.033505508549673
$ws.Range("L20").Value = 1.029696181034062
$ws.Range("N20").Value = 1.014738983769194
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024962257296053
$ws.Range("D21").Value = 1.029417683273709
$ws.Range("E21").Value = 1.025475606678611
$ws.Range("I21").Value = 1.032911874882057
$ws.Range("J21").Value = 1.031628531669355
$ws.Range("K21").Value = 1.033022491652762
$ws.Range("L21").Value = 1.029095357848268
$ws.Range("N21").Value = 1.01455108799268
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024356389674189
$ws.Range("D22").Value = 1.028982022516613
$ws.Range("E22").Value = 1.024965479099662
$ws.Range("I22").Value = 1.032782770578668
$ws.Range("J22").Value = 1.031274454333288
$ws.Range("K22").Value = 1.032718469777914
$ws.Range("L22").Value = 1.028717754192158
$ws.Range("N22").Value = 1.014432859631366
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024677521266839
$ws.Range("D23").Value = 1.029212933918602
$ws.Range("E23").Value = 1.025235823179305
$ws.Range("I23").Value = 1.032851307843285
$ws.Range("J23").Value = 1.031462172688433
$ws.Range("K23").Value = 1.032879668401987
$ws.Range("L23").Value = 1.028917913633439
$ws.Range("N23").Value = 1.014495543179203
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025942893202999
$ws.Range("D24").Value = 1.030122901228999
$ws.Range("E24").Value = 1.026301999702164
$ws.Range("I24").Value = 1.033118985934266
$ws.Range("J24").Value = 1.032200856736259
$ws.Range("K24").Value = 1.033513600709106
$ws.Range("L24").Value = 1.029706256488583
$ws.Range("N24").Value = 1.014742132293209
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027413523609534
$ws.Range("D25").Value = 1.031180648790175
$ws.Range("E25").Value = 1.027543019163745
$ws.Range("I25").Value = 1.033425140172777
$ws.Range("J25").Value = 1.033057296029411
$ws.Range("K25").Value = 1.034247766500302
$ws.Range("L25").Value = 1.030621725916716
$ws.Range("N25").Value = 1.015027874323378
